$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per latest crypto data refresh.
# Price values are forced to Text format so that values such as "1.00" or
# "14.50" retain their original textual representation instead of being
# reinterpreted as numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.161.59"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.280.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "155.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +15,447.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "95.34"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.69%  "

$ws.Range("E8").Value = "  -0.18%  "

$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.494"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "35.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0804"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.71"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.633.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.271.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.799"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.074.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("E31").Value = "  +1.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.44%  "

$ws.Range("E35").Value = "  +0.30%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.58%  "

$ws.Range("E37").Value = "  +4.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.29%  "

$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("E40").Value = "  -0.38%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("E42").Value = "  +6.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.012.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.78%  "

$ws.Range("E46").Value = "  +1.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.52%  "

$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("E51").Value = "  -0.63%  "
